# B6-PowerPoint.pptx edit: re-style the three data tables (slides 14-16)
# from the plain "Table_0" style to the built-in PowerPoint table style
# {5303C78C-D13B-43BE-A36F-4BDC19266245}.
#
# PowerPoint's Table object does not allow TableStyleId to be written as a
# plain property (it throws "Table styles cannot be assigned through a
# property ... call Table.ApplyStyle(...) instead"), so ApplyStyle is used.

$p = $ppt.ActivePresentation

$targetStyleId = "{5303C78C-D13B-43BE-A36F-4BDC19266245}"

for ($slideIdx = 1; $slideIdx -le $p.Slides.Count; $slideIdx++) {
    $slide = $p.Slides.Item($slideIdx)

    for ($shapeIdx = 1; $shapeIdx -le $slide.Shapes.Count; $shapeIdx++) {
        $shape = $slide.Shapes.Item($shapeIdx)

        if ($shape.HasTable) {
            $table = $shape.Table
            $table.ApplyStyle($targetStyleId)
        }
    }
}
